# multi-BU support for supplier list: edge labels, supplier_id dedup, business_unit column
#
# This script rewrites the "Supplier List" header row (row 4) to insert two new
# columns - "supplier_id" (after supplier_name) and "business_unit" (after
# parent_supplier) - shifts the remaining headers right, resizes every column
# to its new target width, and updates the autoFilter / defined name /
# data-validation ranges that referenced the old (18-column) layout so that
# they point at the new (20-column) layout instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 4): write the new, shifted set of header labels.
#    supplier_id is new at column B; business_unit is new at column F.
#    Everything from C4 (prev B4) through T4 (prev R4) just gets the value
#    of what used to sit one or two columns to the left - this naturally
#    re-uses the existing s="4" header style for every cell that already
#    existed, without disturbing rows 1-2 (which the diff leaves untouched).
# ---------------------------------------------------------------------------
$headers = @(
    "supplier_name",
    "supplier_id",
    "jurisdiction",
    "tier",
    "parent_supplier",
    "business_unit",
    "commodity",
    "valid_from",
    "annual_value",
    "value_currency",
    "contract_ref",
    "lei",
    "duns",
    "vat",
    "vat_country",
    "internal_id",
    "risk_tier",
    "kraljic_quadrant",
    "approval_status",
    "notes"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $headers[$i]
}

# The header row grew from 18 to 20 columns, so the brand-new cells
# (supplier_id at B4, business_unit at F4, and the two cells that now exist
# past the old last column R4, namely S4/T4) start out with the default
# style. Paste the header format across the whole row so every header cell
# - old and new - ends up sharing the same header style (s="4") as before,
# reusing the existing style index instead of minting a near-duplicate one.
$ws.Range("A4").Copy()
$ws.Range("A4:T4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Column widths: set every column (1-20) to its target stored width.
#    Excel's ColumnWidth property is in "characters" and is offset from the
#    stored OOXML width by the standard 5px/6 (~0.8333) padding constant, so
#    subtract that constant to land exactly on the desired stored width.
# ---------------------------------------------------------------------------
$padding = 5 / 6
$targetWidths = @(30, 14, 14, 8, 30, 18, 20, 14, 14, 14, 16, 24, 14, 20, 14, 16, 12, 18, 16, 30)

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $ws.Columns($i + 1).ColumnWidth = $targetWidths[$i] - $padding
}

# ---------------------------------------------------------------------------
# 3. autoFilter: re-point it at the new header range A4:T4.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
[void]$ws.Range("A4:T4").AutoFilter()

# ---------------------------------------------------------------------------
# 4. Defined name _xlnm._FilterDatabase: move from $A$4:$R$4 to $A$4:$T$4.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Supplier List'!`$A`$4:`$T`$4"
    }
}

# ---------------------------------------------------------------------------
# 5. Data validations: move each rule to its new column, preserving the
#    original list values, prompts and alert behaviour.
#    tier:             C5:C10000 -> D5:D10000
#    risk_tier:        O5:O10000 -> Q5:Q10000
#    kraljic_quadrant: P5:P10000 -> R5:R10000
#    approval_status:  Q5:Q10000 -> S5:S10000
# ---------------------------------------------------------------------------
$ws.Range("C5:C10000").Validation.Delete()
$ws.Range("O5:O10000").Validation.Delete()
$ws.Range("P5:P10000").Validation.Delete()
$ws.Range("Q5:Q10000").Validation.Delete()

$ws.Range("D5:D10000").Validation.Add(3, 1, 1, '"1,2,3"')
$ws.Range("D5:D10000").Validation.IgnoreBlank = $true
$ws.Range("D5:D10000").Validation.InCellDropdown = $true
$ws.Range("D5:D10000").Validation.ShowInput = $true
$ws.Range("D5:D10000").Validation.ShowError = $true
$ws.Range("D5:D10000").Validation.InputTitle = "Tier"
$ws.Range("D5:D10000").Validation.InputMessage = "Supply-chain tier: 1 = direct, 2 = sub-supplier, 3 = sub-sub-supplier"

$ws.Range("Q5:Q10000").Validation.Add(3, 1, 1, '"critical,high,medium,low"')
$ws.Range("Q5:Q10000").Validation.IgnoreBlank = $true
$ws.Range("Q5:Q10000").Validation.InCellDropdown = $true
$ws.Range("Q5:Q10000").Validation.ShowInput = $true
$ws.Range("Q5:Q10000").Validation.ShowError = $false
$ws.Range("Q5:Q10000").Validation.InputTitle = "Risk Tier"
$ws.Range("Q5:Q10000").Validation.InputMessage = "General risk classification"

$ws.Range("R5:R10000").Validation.Add(3, 1, 1, '"strategic,leverage,bottleneck,non-critical"')
$ws.Range("R5:R10000").Validation.IgnoreBlank = $true
$ws.Range("R5:R10000").Validation.InCellDropdown = $true
$ws.Range("R5:R10000").Validation.ShowInput = $true
$ws.Range("R5:R10000").Validation.ShowError = $false
$ws.Range("R5:R10000").Validation.InputTitle = "Kraljic Quadrant"
$ws.Range("R5:R10000").Validation.InputMessage = "Kraljic portfolio classification"

$ws.Range("S5:S10000").Validation.Add(3, 1, 1, '"approved,conditional,pending,blocked,phase-out"')
$ws.Range("S5:S10000").Validation.IgnoreBlank = $true
$ws.Range("S5:S10000").Validation.InCellDropdown = $true
$ws.Range("S5:S10000").Validation.ShowInput = $true
$ws.Range("S5:S10000").Validation.ShowError = $false
$ws.Range("S5:S10000").Validation.InputTitle = "Approval Status"
$ws.Range("S5:S10000").Validation.InputMessage = "Supplier approval status"

Write-Host "Applied multi-BU supplier list header changes."
